$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had the "Split ^^^" header row at row 6, followed by
# the two data rows (1's then 2's) at rows 7-8. Running the test again with
# the rows re-numbered 2 & 5 means two fresh blank rows get inserted above
# the "Split ^^^" row, pushing it (and the data rows beneath it) down by two
# rows: "Split ^^^" -> row 8, the 1's -> row 9, the 2's -> row 10.
$ws.Rows("6:7").Insert()

# Excel's row insert copies the formatting of the row above into the new
# rows; set the new cells back to blank text (leading apostrophe = empty
# text entry) and strip that inherited formatting so the new rows 6-7 are
# plain/unstyled blank cells, matching the rest of the freshly inserted area.
$ws.Range("A6:F7").Value = "'"
$ws.Rows("6:7").ClearFormats()
